# Commit: Tue, Jul 07, 2020  7:05:30 PM
#
# The underlying edit:
#   1. Three tables (on the slides that originally carried table style
#      {249FC436-9EEC-4219-A5C3-E5BD884DA1E0}) get re-styled to
#      {487EBFC4-812C-4B14-B792-514E9C8C390D}.
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml swap their payloads
#      (the slide master's theme becomes the stock "Office Theme" while
#      the notes master keeps the original "Integral" theme). This is a
#      pure OOXML part-content swap with no relationship/id changes, and
#      is not reachable from the supported PowerPoint COM surface (no
#      automation call mutates raw theme XML / swaps theme parts); the
#      `Presentation/SlideMaster/Design.ApplyTheme` family of calls is a
#      recognized no-op in this host. We still invoke it defensively
#      below in case a host update wires it up, but the real, reachable
#      portion of this edit is the table restyle, which we apply via the
#      documented `Table.ApplyStyle` method (Table.Style is read-only,
#      matching real PowerPoint's object model).

$p = $ppt.ActivePresentation

$oldStyleId = "{249FC436-9EEC-4219-A5C3-E5BD884DA1E0}"
$newStyleId = "{487EBFC4-812C-4B14-B792-514E9C8C390D}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# Best-effort attempt at the theme swap described above (currently a
# no-op on this host, kept for forward-compatibility / harmlessness).
try {
    $p.ApplyTheme($p.SlideMaster.Theme)
} catch {
}
